# Update the timestamp portion of the test email addresses in the
# "UsuariosRegistro" sheet from 20251109_011412 to 20251109_012452.
# Because these email strings are also referenced (as shared strings)
# from the "LoginData" sheet, updating the source cells will keep both
# sheets consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UsuariosRegistro")

$ws.Range("C2").Value = "juan.perez+20251109_012452@test.com"
$ws.Range("C3").Value = "maria.gonzalez+20251109_012452@test.com"
$ws.Range("C4").Value = "carlos.rodriguez+20251109_012452@test.com"
$ws.Range("C5").Value = "ana.martinez+20251109_012452@test.com"
$ws.Range("C6").Value = "luis.garcia+20251109_012452@test.com"

# The "LoginData" sheet references the same juan.perez / maria.gonzalez
# addresses; make sure they reflect the updated timestamp as well.
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Range("A2").Value = "juan.perez+20251109_012452@test.com"
$wsLogin.Range("A3").Value = "maria.gonzalez+20251109_012452@test.com"
